$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 3000  # H74: 700 -> 3000
$ws.Cells.Item(74, 9).Value = 0  # I74: 700 -> 0
$ws.Cells.Item(74, 10).Value = 3000  # J74: 0 -> 3000
$ws.Cells.Item(74, 11).Value = 0  # K74: 700 -> 0
$ws.Cells.Item(74, 12).Value = 3000  # L74: 0 -> 3000
$ws.Cells.Item(74, 13).ClearContents()  # M74: 236 -> (removed)
$ws.Cells.Item(74, 14).Value = -4872  # N74: None -> -4872
$ws.Cells.Item(77, 8).Value = 3000  # H77: 700 -> 3000
$ws.Cells.Item(77, 9).Value = 0  # I77: 700 -> 0
$ws.Cells.Item(77, 10).Value = 3000  # J77: 0 -> 3000
$ws.Cells.Item(77, 11).Value = 0  # K77: 3500 -> 0
$ws.Cells.Item(77, 12).Value = 15000  # L77: 0 -> 15000
$ws.Cells.Item(77, 13).ClearContents()  # M77: 1180 -> (removed)
$ws.Cells.Item(77, 14).Value = -24360  # N77: None -> -24360
$ws.Cells.Item(88, 8).Value = 1544852.9  # H88: 2059569.9 -> 1544852.9
$ws.Cells.Item(88, 9).Value = 1000  # I88: 0 -> 1000
$ws.Cells.Item(88, 10).Value = 2059470.5  # J88: 2059569.9 -> 2059470.5
$ws.Cells.Item(88, 11).Value = 1000  # K88: 0 -> 1000
$ws.Cells.Item(88, 12).Value = 2059470.5  # L88: 2059569.9 -> 2059470.5
$ws.Cells.Item(88, 13).Value = -594  # M88: None -> -594
$ws.Cells.Item(88, 14).Value = -2060282.5  # N88: -2060381.9 -> -2060282.5
$ws.Cells.Item(91, 8).Value = 1544852.9  # H91: 2059569.9 -> 1544852.9
$ws.Cells.Item(91, 9).Value = 1000  # I91: 0 -> 1000
$ws.Cells.Item(91, 10).Value = 2059470.5  # J91: 2059569.9 -> 2059470.5
$ws.Cells.Item(91, 11).Value = 1000  # K91: 0 -> 1000
$ws.Cells.Item(91, 12).Value = 2059470.5  # L91: 2059569.9 -> 2059470.5
$ws.Cells.Item(91, 13).Value = 404  # M91: None -> 404
$ws.Cells.Item(91, 14).Value = -2062278.5  # N91: -2062377.9 -> -2062278.5
$ws.Cells.Item(116, 8).Value = 2411.5862  # H116: 2387.2 -> 2411.5862
$ws.Cells.Item(116, 9).Value = 2324.0527  # I116: 2291.85 -> 2324.0527
$ws.Cells.Item(116, 11).Value = 2324.0527  # K116: 2291.85 -> 2324.0527
$ws.Cells.Item(116, 13).Value = 1117.9473  # M116: 1150.15 -> 1117.9473
$ws.Cells.Item(132, 8).Value = 6809263.5  # H132: 7099005 -> 6809263.5
$ws.Cells.Item(132, 9).Value = 9528356  # I132: 9808600 -> 9528356
$ws.Cells.Item(132, 10).Value = 11533.786  # J132: 12371 -> 11533.786
$ws.Cells.Item(132, 11).Value = 28585068  # K132: 29425800 -> 28585068
$ws.Cells.Item(132, 12).Value = 34601.358  # L132: 37113 -> 34601.358
$ws.Cells.Item(132, 13).Value = -28582538  # M132: -29423270 -> -28582538
$ws.Cells.Item(132, 14).Value = -39661.358  # N132: -42173 -> -39661.358
$ws.Cells.Item(138, 8).Value = 1535.7374  # H138: 1553.64 -> 1535.7374
$ws.Cells.Item(138, 9).Value = 1092.1538  # I138: 1099.8334 -> 1092.1538
$ws.Cells.Item(138, 10).Value = 1602.7906  # J138: 1615.5227 -> 1602.7906
$ws.Cells.Item(138, 11).Value = 3276.4614  # K138: 3299.5002 -> 3276.4614
$ws.Cells.Item(138, 12).Value = 4808.3718  # L138: 4846.5681 -> 4808.3718
$ws.Cells.Item(138, 13).Value = 1863.5386  # M138: 1840.4998 -> 1863.5386
$ws.Cells.Item(138, 14).Value = -15088.3718  # N138: -15126.5681 -> -15088.3718

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4009.554  # H32: 3934.7163 -> 4009.554
$ws.Cells.Item(32, 9).Value = 3468.8667  # I32: 3405.4355 -> 3468.8667
$ws.Cells.Item(32, 11).Value = 3468.8667  # K32: 3405.4355 -> 3468.8667
$ws.Cells.Item(32, 13).Value = -3181.8667  # M32: -3118.4355 -> -3181.8667
$ws.Cells.Item(61, 8).Value = 1704.5769  # H61: 1719.96 -> 1704.5769
$ws.Cells.Item(61, 9).Value = 1856.1765  # I61: 1834.4117 -> 1856.1765
$ws.Cells.Item(61, 10).Value = 1418.2222  # J61: 1476.75 -> 1418.2222
$ws.Cells.Item(61, 11).Value = 1856.1765  # K61: 1834.4117 -> 1856.1765
$ws.Cells.Item(61, 12).Value = 1418.2222  # L61: 1476.75 -> 1418.2222
$ws.Cells.Item(61, 13).Value = -1644.1765  # M61: -1622.4117 -> -1644.1765
$ws.Cells.Item(61, 14).Value = -1842.2222  # N61: -1900.75 -> -1842.2222
$ws.Cells.Item(102, 8).Value = 55556424  # H102: 55556724 -> 55556424
$ws.Cells.Item(102, 9).Value = 55556424  # I102: 55556724 -> 55556424
$ws.Cells.Item(102, 11).Value = 55556424  # K102: 55556724 -> 55556424
$ws.Cells.Item(102, 13).Value = -55554802  # M102: -55555102 -> -55554802
$ws.Cells.Item(122, 8).Value = 1356.5883  # H122: 1266.6666 -> 1356.5883
$ws.Cells.Item(122, 9).Value = 1257.5385  # I122: 1500 -> 1257.5385
$ws.Cells.Item(122, 10).Value = 1678.5  # J122: 566.6667 -> 1678.5
$ws.Cells.Item(122, 11).Value = 3772.6155  # K122: 4500 -> 3772.6155
$ws.Cells.Item(122, 12).Value = 5035.5  # L122: 1700.0001 -> 5035.5
$ws.Cells.Item(122, 13).Value = -1322.6155  # M122: -2050 -> -1322.6155
$ws.Cells.Item(122, 14).Value = -9935.5  # N122: -6600.0001 -> -9935.5
$ws.Cells.Item(132, 8).Value = 1153.7188  # H132: 1191.3226 -> 1153.7188
$ws.Cells.Item(132, 9).Value = 871.9245  # I132: 918.4400000000001 -> 871.9245
$ws.Cells.Item(132, 10).Value = 2511.4546  # J132: 2328.3333 -> 2511.4546
$ws.Cells.Item(132, 11).Value = 2615.7735  # K132: 2755.32 -> 2615.7735
$ws.Cells.Item(132, 12).Value = 7534.3638  # L132: 6984.999899999999 -> 7534.3638
$ws.Cells.Item(132, 13).Value = -85.77349999999979  # M132: -225.3200000000002 -> -85.77349999999979
$ws.Cells.Item(132, 14).Value = -12594.3638  # N132: -12044.9999 -> -12594.3638
$ws.Cells.Item(136, 8).Value = 1704.5769  # H136: 1719.96 -> 1704.5769
$ws.Cells.Item(136, 9).Value = 1856.1765  # I136: 1834.4117 -> 1856.1765
$ws.Cells.Item(136, 10).Value = 1418.2222  # J136: 1476.75 -> 1418.2222
$ws.Cells.Item(136, 11).Value = 5568.529500000001  # K136: 5503.2351 -> 5568.529500000001
$ws.Cells.Item(136, 12).Value = 4254.6666  # L136: 4430.25 -> 4254.6666
$ws.Cells.Item(136, 13).Value = -3018.529500000001  # M136: -2953.2351 -> -3018.529500000001
$ws.Cells.Item(136, 14).Value = -9354.6666  # N136: -9530.25 -> -9354.6666

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2097.5557  # H20: 2344.3333 -> 2097.5557
$ws.Cells.Item(20, 9).Value = 2239.818  # I20: 2538.7778 -> 2239.818
$ws.Cells.Item(20, 10).Value = 1874  # J20: 2052.6667 -> 1874
$ws.Cells.Item(20, 11).Value = 2239.818  # K20: 2538.7778 -> 2239.818
$ws.Cells.Item(20, 12).Value = 1874  # L20: 2052.6667 -> 1874
$ws.Cells.Item(20, 13).Value = -1992.818  # M20: -2291.7778 -> -1992.818
$ws.Cells.Item(20, 14).Value = -2368  # N20: -2546.6667 -> -2368
$ws.Cells.Item(86, 8).Value = 3882.5293  # H86: 4205.8 -> 3882.5293
$ws.Cells.Item(86, 9).Value = 3882.5293  # I86: 4205.8 -> 3882.5293
$ws.Cells.Item(86, 11).Value = 3882.5293  # K86: 4205.8 -> 3882.5293
$ws.Cells.Item(86, 13).Value = -2759.5293  # M86: -3082.8 -> -2759.5293
$ws.Cells.Item(89, 8).Value = 3882.5293  # H89: 4205.8 -> 3882.5293
$ws.Cells.Item(89, 9).Value = 3882.5293  # I89: 4205.8 -> 3882.5293
$ws.Cells.Item(89, 11).Value = 19412.6465  # K89: 21029 -> 19412.6465
$ws.Cells.Item(89, 13).Value = -13796.6465  # M89: -15413 -> -13796.6465
$ws.Cells.Item(94, 8).Value = 15625887  # H94: 20834228 -> 15625887
$ws.Cells.Item(94, 9).Value = 16667465  # I94: 20834228 -> 16667465
$ws.Cells.Item(94, 10).Value = 2222  # J94: 0 -> 2222
$ws.Cells.Item(94, 11).Value = 16667465  # K94: 20834228 -> 16667465
$ws.Cells.Item(94, 12).Value = 2222  # L94: 0 -> 2222
$ws.Cells.Item(94, 13).Value = -16667014  # M94: -20833777 -> -16667014
$ws.Cells.Item(94, 14).Value = -3124  # N94: None -> -3124
$ws.Cells.Item(134, 8).Value = 4259.3335  # H134: 5071.5625 -> 4259.3335
$ws.Cells.Item(134, 9).Value = 967.74194  # I134: 1079.6 -> 967.74194
$ws.Cells.Item(134, 10).Value = 17014.25  # J134: 19328.572 -> 17014.25
$ws.Cells.Item(134, 11).Value = 2903.22582  # K134: 3238.8 -> 2903.22582
$ws.Cells.Item(134, 12).Value = 51042.75  # L134: 57985.716 -> 51042.75
$ws.Cells.Item(134, 13).Value = -368.2258200000001  # M134: -703.7999999999997 -> -368.2258200000001
$ws.Cells.Item(134, 14).Value = -56112.75  # N134: -63055.716 -> -56112.75

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 10426.667  # H23: 10000 -> 10426.667
$ws.Cells.Item(23, 9).Value = 9900  # I23: 0 -> 9900
$ws.Cells.Item(23, 10).Value = 10492.5  # J23: 10000 -> 10492.5
$ws.Cells.Item(23, 11).Value = 9900  # K23: 0 -> 9900
$ws.Cells.Item(23, 12).Value = 10492.5  # L23: 10000 -> 10492.5
$ws.Cells.Item(23, 13).Value = -9660  # M23: None -> -9660
$ws.Cells.Item(23, 14).Value = -10972.5  # N23: -10480 -> -10972.5
$ws.Cells.Item(27, 8).Value = 10426.667  # H27: 10000 -> 10426.667
$ws.Cells.Item(27, 9).Value = 9900  # I27: 0 -> 9900
$ws.Cells.Item(27, 10).Value = 10492.5  # J27: 10000 -> 10492.5
$ws.Cells.Item(27, 11).Value = 9900  # K27: 0 -> 9900
$ws.Cells.Item(27, 12).Value = 10492.5  # L27: 10000 -> 10492.5
$ws.Cells.Item(27, 13).Value = -9708  # M27: None -> -9708
$ws.Cells.Item(27, 14).Value = -10876.5  # N27: -10384 -> -10876.5
$ws.Cells.Item(41, 8).Value = 9850  # H41: 13350 -> 9850
$ws.Cells.Item(41, 9).Value = 2550  # I41: 2733.3333 -> 2550
$ws.Cells.Item(41, 10).Value = 24450  # J41: 23966.666 -> 24450
$ws.Cells.Item(41, 11).Value = 2550  # K41: 2733.3333 -> 2550
$ws.Cells.Item(41, 12).Value = 24450  # L41: 23966.666 -> 24450
$ws.Cells.Item(41, 13).Value = -2122  # M41: -2305.3333 -> -2122
$ws.Cells.Item(41, 14).Value = -25306  # N41: -24822.666 -> -25306
$ws.Cells.Item(58, 8).Value = 567.9056399999999  # H58: 613.40424 -> 567.9056399999999
$ws.Cells.Item(58, 9).Value = 561.40625  # I58: 671.2917 -> 561.40625
$ws.Cells.Item(58, 10).Value = 577.8095  # J58: 553 -> 577.8095
$ws.Cells.Item(58, 11).Value = 561.40625  # K58: 671.2917 -> 561.40625
$ws.Cells.Item(58, 12).Value = 577.8095  # L58: 553 -> 577.8095
$ws.Cells.Item(58, 13).Value = -358.40625  # M58: -468.2917 -> -358.40625
$ws.Cells.Item(58, 14).Value = -983.8095  # N58: -959 -> -983.8095
$ws.Cells.Item(114, 8).Value = 29490  # H114: 26660 -> 29490
$ws.Cells.Item(114, 9).Value = 0  # I114: 21000 -> 0
$ws.Cells.Item(114, 11).Value = 0  # K114: 21000 -> 0
$ws.Cells.Item(114, 13).ClearContents()  # M114: -16661 -> (removed)
$ws.Cells.Item(132, 8).Value = 1426.8776  # H132: 1664.9286 -> 1426.8776
$ws.Cells.Item(132, 9).Value = 1075.0975  # I132: 1296.7354 -> 1075.0975
$ws.Cells.Item(132, 11).Value = 3225.2925  # K132: 3890.2062 -> 3225.2925
$ws.Cells.Item(132, 13).Value = -695.2925000000005  # M132: -1360.2062 -> -695.2925000000005
$ws.Cells.Item(134, 8).Value = 658.39026  # H134: 604.9792 -> 658.39026
$ws.Cells.Item(134, 9).Value = 631.06665  # I134: 575.75 -> 631.06665
$ws.Cells.Item(134, 10).Value = 732.9091  # J134: 692.6667 -> 732.9091
$ws.Cells.Item(134, 11).Value = 1893.19995  # K134: 1727.25 -> 1893.19995
$ws.Cells.Item(134, 12).Value = 2198.7273  # L134: 2078.0001 -> 2198.7273
$ws.Cells.Item(134, 13).Value = 641.8000500000001  # M134: 807.75 -> 641.8000500000001
$ws.Cells.Item(134, 14).Value = -7268.7273  # N134: -7148.0001 -> -7268.7273
$ws.Cells.Item(136, 8).Value = 567.9056399999999  # H136: 613.40424 -> 567.9056399999999
$ws.Cells.Item(136, 9).Value = 561.40625  # I136: 671.2917 -> 561.40625
$ws.Cells.Item(136, 10).Value = 577.8095  # J136: 553 -> 577.8095
$ws.Cells.Item(136, 11).Value = 1684.21875  # K136: 2013.8751 -> 1684.21875
$ws.Cells.Item(136, 12).Value = 1733.4285  # L136: 1659 -> 1733.4285
$ws.Cells.Item(136, 13).Value = 865.78125  # M136: 536.1249 -> 865.78125
$ws.Cells.Item(136, 14).Value = -6833.4285  # N136: -6759 -> -6833.4285
$ws.Cells.Item(141, 8).Value = 27740  # H141: 28080 -> 27740
$ws.Cells.Item(141, 10).Value = 27740  # J141: 28080 -> 27740
$ws.Cells.Item(141, 12).Value = 27740  # L141: 28080 -> 27740
$ws.Cells.Item(141, 14).Value = -38100  # N141: -38440 -> -38100

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 163223.11  # H11: 163212.62 -> 163223.11
$ws.Cells.Item(11, 9).Value = 182378.77  # I11: 182390.6 -> 182378.77
$ws.Cells.Item(11, 10).Value = 400  # J11: 200 -> 400
$ws.Cells.Item(11, 11).Value = 547136.3099999999  # K11: 547171.8 -> 547136.3099999999
$ws.Cells.Item(11, 12).Value = 1200  # L11: 600 -> 1200
$ws.Cells.Item(11, 13).Value = -546996.3099999999  # M11: -547031.8 -> -546996.3099999999
$ws.Cells.Item(11, 14).Value = -1480  # N11: -880 -> -1480
$ws.Cells.Item(86, 8).Value = 570  # H86: 446.66666 -> 570
$ws.Cells.Item(86, 9).Value = 0  # I86: 200 -> 0
$ws.Cells.Item(86, 11).Value = 0  # K86: 600 -> 0
$ws.Cells.Item(86, 13).ClearContents()  # M86: 586 -> (removed)
$ws.Cells.Item(89, 8).Value = 570  # H89: 446.66666 -> 570
$ws.Cells.Item(89, 9).Value = 0  # I89: 200 -> 0
$ws.Cells.Item(89, 11).Value = 0  # K89: 1800 -> 0
$ws.Cells.Item(89, 13).ClearContents()  # M89: 4128 -> (removed)
$ws.Cells.Item(107, 8).Value = 7092.7334  # H107: 6712.0625 -> 7092.7334
$ws.Cells.Item(107, 10).Value = 10440.5  # J107: 9582.454 -> 10440.5
$ws.Cells.Item(107, 12).Value = 31321.5  # L107: 28747.362 -> 31321.5
$ws.Cells.Item(107, 14).Value = -35161.5  # N107: -32587.362 -> -35161.5
$ws.Cells.Item(131, 8).Value = 35715720  # H131: 25642326 -> 35715720
$ws.Cells.Item(131, 10).Value = 2206  # J131: 1717.2858 -> 2206
$ws.Cells.Item(131, 12).Value = 6618  # L131: 5151.857400000001 -> 6618
$ws.Cells.Item(131, 14).Value = -16698  # N131: -15231.8574 -> -16698
$ws.Cells.Item(132, 8).Value = 3500  # H132: 2750 -> 3500
$ws.Cells.Item(132, 10).Value = 0  # J132: 2000 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 18000 -> 0
$ws.Cells.Item(132, 14).ClearContents()  # N132: -23060 -> (removed)

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 75003720  # H70: 64289690 -> 75003720
$ws.Cells.Item(70, 9).Value = 62504176  # I70: 50004440 -> 62504176
$ws.Cells.Item(70, 11).Value = 62504176  # K70: 50004440 -> 62504176
$ws.Cells.Item(70, 13).Value = -62503906  # M70: -50004170 -> -62503906
$ws.Cells.Item(73, 8).Value = 75003720  # H73: 64289690 -> 75003720
$ws.Cells.Item(73, 9).Value = 62504176  # I73: 50004440 -> 62504176
$ws.Cells.Item(73, 11).Value = 62504176  # K73: 50004440 -> 62504176
$ws.Cells.Item(73, 13).Value = -62503240  # M73: -50003504 -> -62503240
$ws.Cells.Item(80, 8).Value = 3228.5715  # H80: 3466.6667 -> 3228.5715
$ws.Cells.Item(83, 8).Value = 3228.5715  # H83: 3466.6667 -> 3228.5715
$ws.Cells.Item(132, 8).Value = 2026.1765  # H132: 2483.3 -> 2026.1765
$ws.Cells.Item(132, 9).Value = 1570.8462  # I132: 1831.7142 -> 1570.8462
$ws.Cells.Item(132, 10).Value = 3506  # J132: 4003.6667 -> 3506
$ws.Cells.Item(132, 11).Value = 4712.5386  # K132: 5495.142599999999 -> 4712.5386
$ws.Cells.Item(132, 12).Value = 10518  # L132: 12011.0001 -> 10518
$ws.Cells.Item(132, 13).Value = -2182.5386  # M132: -2965.142599999999 -> -2182.5386
$ws.Cells.Item(132, 14).Value = -15578  # N132: -17071.0001 -> -15578

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 933.4286  # H16: 875.125 -> 933.4286
$ws.Cells.Item(16, 9).Value = 813.2727  # I16: 782.5 -> 813.2727
$ws.Cells.Item(16, 10).Value = 1374  # J16: 1153 -> 1374
$ws.Cells.Item(16, 11).Value = 813.2727  # K16: 782.5 -> 813.2727
$ws.Cells.Item(16, 12).Value = 1374  # L16: 1153 -> 1374
$ws.Cells.Item(16, 13).Value = -643.2727  # M16: -612.5 -> -643.2727
$ws.Cells.Item(16, 14).Value = -1714  # N16: -1493 -> -1714
$ws.Cells.Item(22, 8).Value = 1040.2  # H22: 1667.6666 -> 1040.2
$ws.Cells.Item(22, 9).Value = 800.5  # I22: 1001 -> 800.5
$ws.Cells.Item(22, 10).Value = 1200  # J22: 2001 -> 1200
$ws.Cells.Item(22, 11).Value = 800.5  # K22: 1001 -> 800.5
$ws.Cells.Item(22, 12).Value = 1200  # L22: 2001 -> 1200
$ws.Cells.Item(22, 13).Value = -505.5  # M22: -706 -> -505.5
$ws.Cells.Item(22, 14).Value = -1790  # N22: -2591 -> -1790
$ws.Cells.Item(27, 8).Value = 1040.2  # H27: 1667.6666 -> 1040.2
$ws.Cells.Item(27, 9).Value = 800.5  # I27: 1001 -> 800.5
$ws.Cells.Item(27, 10).Value = 1200  # J27: 2001 -> 1200
$ws.Cells.Item(27, 11).Value = 800.5  # K27: 1001 -> 800.5
$ws.Cells.Item(27, 12).Value = 1200  # L27: 2001 -> 1200
$ws.Cells.Item(27, 13).Value = -693.5  # M27: -894 -> -693.5
$ws.Cells.Item(27, 14).Value = -1414  # N27: -2215 -> -1414
$ws.Cells.Item(46, 8).Value = 1143.75  # H46: 1350 -> 1143.75
$ws.Cells.Item(46, 9).Value = 917.5  # I46: 1000 -> 917.5
$ws.Cells.Item(46, 10).Value = 1370  # J46: 1700 -> 1370
$ws.Cells.Item(46, 11).Value = 917.5  # K46: 1000 -> 917.5
$ws.Cells.Item(46, 12).Value = 1370  # L46: 1700 -> 1370
$ws.Cells.Item(46, 13).Value = -729.5  # M46: -812 -> -729.5
$ws.Cells.Item(46, 14).Value = -1746  # N46: -2076 -> -1746
$ws.Cells.Item(82, 8).Value = 1601.4615  # H82: 1657.9615 -> 1601.4615
$ws.Cells.Item(82, 10).Value = 1569.0667  # J82: 1667 -> 1569.0667
$ws.Cells.Item(82, 12).Value = 1569.0667  # L82: 1667 -> 1569.0667
$ws.Cells.Item(82, 14).Value = -2291.0667  # N82: -2389 -> -2291.0667
$ws.Cells.Item(85, 8).Value = 1601.4615  # H85: 1657.9615 -> 1601.4615
$ws.Cells.Item(85, 10).Value = 1569.0667  # J85: 1667 -> 1569.0667
$ws.Cells.Item(85, 12).Value = 1569.0667  # L85: 1667 -> 1569.0667
$ws.Cells.Item(85, 14).Value = -4065.0667  # N85: -4163 -> -4065.0667
$ws.Cells.Item(122, 8).Value = 31252738  # H122: 35717130 -> 31252738
$ws.Cells.Item(122, 9).Value = 83335660  # I122: 125002500 -> 83335660
$ws.Cells.Item(122, 11).Value = 250006980  # K122: 375007500 -> 250006980
$ws.Cells.Item(122, 13).Value = -250004530  # M122: -375005050 -> -250004530
$ws.Cells.Item(132, 8).Value = 22150.898  # H132: 22108.041 -> 22150.898
$ws.Cells.Item(132, 9).Value = 1126.4231  # I132: 1118.037 -> 1126.4231
$ws.Cells.Item(132, 10).Value = 45917.695  # J132: 47868.5 -> 45917.695
$ws.Cells.Item(132, 11).Value = 3379.2693  # K132: 3354.111 -> 3379.2693
$ws.Cells.Item(132, 12).Value = 137753.085  # L132: 143605.5 -> 137753.085
$ws.Cells.Item(132, 13).Value = -849.2692999999999  # M132: -824.1109999999999 -> -849.2692999999999
$ws.Cells.Item(132, 14).Value = -142813.085  # N132: -148665.5 -> -142813.085
$ws.Cells.Item(136, 8).Value = 1143.9429  # H136: 1328.4445 -> 1143.9429
$ws.Cells.Item(136, 9).Value = 931.0741  # I136: 1122.3334 -> 931.0741
$ws.Cells.Item(136, 10).Value = 1862.375  # J136: 2049.8333 -> 1862.375
$ws.Cells.Item(136, 11).Value = 2793.2223  # K136: 3367.0002 -> 2793.2223
$ws.Cells.Item(136, 12).Value = 5587.125  # L136: 6149.499899999999 -> 5587.125
$ws.Cells.Item(136, 13).Value = -243.2223000000004  # M136: -817.0001999999999 -> -243.2223000000004
$ws.Cells.Item(136, 14).Value = -10687.125  # N136: -11249.4999 -> -10687.125

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 445.35715  # H113: 486.25 -> 445.35715
$ws.Cells.Item(113, 9).Value = 333.5  # I113: 359.44446 -> 333.5
$ws.Cells.Item(113, 10).Value = 725  # J113: 866.6667 -> 725
$ws.Cells.Item(113, 11).Value = 1000.5  # K113: 1078.33338 -> 1000.5
$ws.Cells.Item(113, 12).Value = 2175  # L113: 2600.0001 -> 2175
$ws.Cells.Item(113, 13).Value = 1169.5  # M113: 1091.66662 -> 1169.5
$ws.Cells.Item(113, 14).Value = -6515  # N113: -6940.0001 -> -6515
$ws.Cells.Item(122, 8).Value = 17858836  # H122: 16668314 -> 17858836
$ws.Cells.Item(122, 9).Value = 19232516  # I122: 17858836 -> 19232516
$ws.Cells.Item(122, 11).Value = 57697548  # K122: 53576508 -> 57697548
$ws.Cells.Item(122, 13).Value = -57695098  # M122: -53574058 -> -57695098
$ws.Cells.Item(132, 8).Value = 2107.3044  # H132: 1829.6852 -> 2107.3044
$ws.Cells.Item(132, 9).Value = 1975.8636  # I132: 1707.7885 -> 1975.8636
$ws.Cells.Item(132, 11).Value = 5927.5908  # K132: 5123.3655 -> 5927.5908
$ws.Cells.Item(132, 13).Value = -3397.5908  # M132: -2593.3655 -> -3397.5908
